$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"
$ws.Range("C2").Value = 2771.04675450926
$ws.Range("C4").Value = 5596.139681459835
$ws.Range("C5").Value = 2870.311589353206
$ws.Range("C6").Value = 1460.056109840828
$ws.Range("C7").Value = 5191.140356354663
$ws.Range("AL7").Value = 1
$ws.Range("C8").Value = 6128.19547247793
$ws.Range("C9").Value = 4729.735976516416
$ws.Range("C10").Value = 5741.405300355145
$ws.Range("C11").Value = 3799.441100542036
$ws.Range("C12").Value = 7397.509860835168
$ws.Range("C14").Value = 2839.92516805933
$ws.Range("C16").Value = 710.2742021758368
$ws.Range("C17").Value = 473.2998774917226
$ws.Range("C18").Value = 5730.354774594881
$ws.Range("C19").Value = 2898.942214704482
$ws.Range("C20").Value = 665.6274194933962
$ws.Range("AL20").Value = 1
$ws.Range("C21").Value = 1503.870423231357
$ws.Range("C22").Value = 5555.389721901988
$ws.Range("AL22").Value = 1
$ws.Range("C23").Value = 6336.709213679884
$ws.Range("C24").Value = 5082.354756663512
$ws.Range("C25").Value = 5814.327715027426
$ws.Range("C26").Value = 4635.517779317834
$ws.Range("C28").Value = 7761.646190572197
$ws.Range("C29").Value = 3690.113267786719
$ws.Range("C30").Value = 2948.84548976845
$ws.Range("C32").Value = 711.1128122770988
$ws.Range("C33").Value = 466.0709276378625
$ws.Range("C34").Value = 5885.254624554112
$ws.Range("C35").Value = 2965.153206179127
$ws.Range("C36").Value = 691.8942672110555
$ws.Range("AL36").Value = 1
$ws.Range("C37").Value = 1577.487171555845
$ws.Range("C38").Value = 5660.517066940175
$ws.Range("AL38").Value = 1
$ws.Range("C39").Value = 6711.616186806423
$ws.Range("C40").Value = 5360.226632400601
$ws.Range("C41").Value = 3587.183047009039
$ws.Range("C42").Value = 7453.823475007535
$ws.Range("C44").Value = 6911.59200404802
$ws.Range("C45").Value = 2999.422762626143
$ws.Range("C46").Value = 6051.685746144485
$ws.Range("C47").Value = 1657.651524528445
$ws.Range("C48").Value = 684.6474015015979
$ws.Range("C49").Value = 495.763971160512
$ws.Range("C50").Value = 2995.45235738661
$ws.Range("C52").Value = 3487.613616731733
$ws.Range("C53").Value = 7179.116970062444
$ws.Range("C55").Value = 7200.731056811853
$ws.Range("C56").Value = 3056.152683606517
$ws.Range("C57").Value = 6203.843262938323
$ws.Range("C58").Value = 1716.389195271215
$ws.Range("C59").Value = 680.3923729568069
$ws.Range("C60").Value = 503.3023574516347
$ws.Range("C61").Value = 3087.12349650562
$ws.Range("C63").Value = 3405.472039138021
$ws.Range("C64").Value = 6978.952586250825
$ws.Range("C66").Value = 7449.08671983612
$ws.Range("C67").Value = 6255.426161047989
$ws.Range("C68").Value = 3137.260298393558
$ws.Range("C69").Value = 3353.623382286602
$ws.Range("C70").Value = 6753.607115829548
$ws.Range("C71").Value = 707.8672001573369
$ws.Range("C72").Value = 3125.07948072635
$ws.Range("C73").Value = 1775.027517189621
$ws.Range("C74").Value = 515.8271637832048
$ws.Range("C77").Value = 7580.275568826287
$ws.Range("C78").Value = 6522.736799041846
$ws.Range("C79").Value = 3210.869677115934
$ws.Range("C80").Value = 3305.422815235401
$ws.Range("C81").Value = 6487.899081675427
$ws.Range("C82").Value = 729.7808175407341
$ws.Range("C83").Value = 3222.05417836739
$ws.Range("C84").Value = 1836.014008604312
$ws.Range("C85").Value = 517.8609592583078
$ws.Range("C88").Value = 7633.969039669125
$ws.Range("C89").Value = 6550.274372976741
$ws.Range("C90").Value = 3242.636921959078
$ws.Range("C91").Value = 3271.088200372761
$ws.Range("C92").Value = 6411.986543373589
$ws.Range("C93").Value = 749.2194349876407
$ws.Range("C94").Value = 3212.81539531051
$ws.Range("C95").Value = 1895.214690888655
$ws.Range("C96").Value = 526.5953412037009
